$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7991798013051437
$ws.Range("C2").Value = 0.08293557258585338
$ws.Range("D2").Value = 0.02163023039290124
$ws.Range("E2").Value = 0.07654431705198661
$ws.Range("F2").Value = 5.863167701030619
$ws.Range("J2").Value = 0.2344677533486674
$ws.Range("K2").Value = 0.8031111687671739
$ws.Range("M2").Value = 0.2915320632618652

$ws.Range("B3").Value = 0.7828399897699967
$ws.Range("C3").Value = 0.08223725791521019
$ws.Range("D3").Value = 0.02002985914090516
$ws.Range("E3").Value = 0.0773216924431761
$ws.Range("F3").Value = 5.662669389073955
$ws.Range("J3").Value = 0.2304004207804908
$ws.Range("K3").Value = 0.7884657395183012
$ws.Range("M3").Value = 0.2902225149972573

$ws.Range("B4").Value = 0.7735445472805793
$ws.Range("C4").Value = 0.08192538117023673
$ws.Range("D4").Value = 0.01907601518124125
$ws.Range("E4").Value = 0.07785054628256383
$ws.Range("F4").Value = 5.540014274568676
$ws.Range("J4").Value = 0.2279588057309354
$ws.Range("K4").Value = 0.7802822408276597
$ws.Range("M4").Value = 0.289688811142355

$ws.Range("B5").Value = 0.7699415970105008
$ws.Range("C5").Value = 0.08182756097062338
$ws.Range("D5").Value = 0.0186944011048098
$ws.Range("E5").Value = 0.07807903105063296
$ws.Range("F5").Value = 5.490141141982747
$ws.Range("J5").Value = 0.2269777854049551
$ws.Range("K5").Value = 0.7771503040802941
$ws.Range("M5").Value = 0.2895392314254863

$ws.Range("B6").Value = 0.7693544938099421
$ws.Range("C6").Value = 0.08181308246076924
$ws.Range("D6").Value = 0.01863145778672504
$ws.Range("E6").Value = 0.07811775463628301
$ws.Range("F6").Value = 5.481866270475905
$ws.Range("J6").Value = 0.226815729021574
$ws.Range("K6").Value = 0.7766424893714543
$ws.Range("M6").Value = 0.2895184932017756

$ws.Range("B7").Value = 0.7734952079948982
$ws.Range("C7").Value = 0.08192394357035937
$ws.Range("D7").Value = 0.01907084011686777
$ws.Range("E7").Value = 0.07785357516701019
$ws.Range("F7").Value = 5.539341228362446
$ws.Range("J7").Value = 0.227945518909948
$ws.Range("K7").Value = 0.7802391815198604
$ws.Range("M7").Value = 0.2896865189981028

$ws.Range("B8").Value = 0.7933925774970589
$ws.Range("C8").Value = 0.08267045868197442
$ws.Range("D8").Value = 0.02107235336563207
$ws.Range("E8").Value = 0.07680166543622668
$ws.Range("F8").Value = 5.793939814886727
$ws.Range("J8").Value = 0.2330537298676845
$ws.Range("K8").Value = 0.7978932417924796
$ws.Range("M8").Value = 0.2910243682619651

$ws.Range("B9").Value = 0.8382822427377619
$ws.Range("C9").Value = 0.08506772311531563
$ws.Range("D9").Value = 0.02523291214784251
$ws.Range("E9").Value = 0.07514736684017898
$ws.Range("F9").Value = 6.296985270176094
$ws.Range("J9").Value = 0.2435165698313853
$ws.Range("K9").Value = 0.8389568695928347
$ws.Range("M9").Value = 0.2957975146957743

$ws.Range("B10").Value = 0.8748771173455907
$ws.Range("C10").Value = 0.0874067071444955
$ws.Range("D10").Value = 0.02844348372066463
$ws.Range("E10").Value = 0.07418042196775865
$ws.Range("F10").Value = 6.669182482728104
$ws.Range("J10").Value = 0.251480908207725
$ws.Range("K10").Value = 0.8730971014587396
$ws.Range("M10").Value = 0.3006223303459308

$ws.Range("B11").Value = 0.8923181208076585
$ws.Range("C11").Value = 0.08859817342950294
$ws.Range("D11").Value = 0.02993967384836083
$ws.Range("E11").Value = 0.07379438130730342
$ws.Range("F11").Value = 6.839140080640959
$ws.Range("J11").Value = 0.2551655899447383
$ws.Range("K11").Value = 0.8895003311529592
$ws.Range("M11").Value = 0.3031052282368378

$ws.Range("B12").Value = 0.8990373010481392
$ws.Range("C12").Value = 0.08906783240894356
$ws.Range("D12").Value = 0.03051154927827326
$ws.Range("E12").Value = 0.07365592996841563
$ws.Range("F12").Value = 6.903595945573443
$ws.Range("J12").Value = 0.2565698372900869
$ws.Range("K12").Value = 0.8958380041115959
$ws.Range("M12").Value = 0.3040869891392646

$ws.Range("B13").Value = 0.8975850980903317
$ws.Range("C13").Value = 0.08896585878451901
$ws.Range("D13").Value = 0.03038814721563199
$ws.Range("E13").Value = 0.07368540407208535
$ws.Range("F13").Value = 6.88970985884589
$ws.Range("J13").Value = 0.2562670088163657
$ws.Range("K13").Value = 0.8944674530290797
$ws.Range("M13").Value = 0.3038736998732716

$ws.Range("B14").Value = 0.8928686111115951
$ws.Range("C14").Value = 0.08863644139344729
$ws.Range("D14").Value = 0.02998661517698054
$ws.Range("E14").Value = 0.07378283586888656
$ws.Range("F14").Value = 6.844440943436041
$ws.Range("J14").Value = 0.2552809386681361
$ws.Range("K14").Value = 0.890019203859481
$ws.Range("M14").Value = 0.3031851650549342

$ws.Range("B15").Value = 0.8899945727056604
$ws.Range("C15").Value = 0.08843707440357207
$ws.Range("D15").Value = 0.02974136055945564
$ws.Range("E15").Value = 0.07384352265485639
$ws.Range("F15").Value = 6.816725130367445
$ws.Range("J15").Value = 0.2546781084032546
$ws.Range("K15").Value = 0.8873109673224349
$ws.Range("M15").Value = 0.3027688310856576

$ws.Range("B16").Value = 0.8737533228087386
$ws.Range("C16").Value = 0.08733141934892785
$ws.Range("D16").Value = 0.02834643669161352
$ws.Range("E16").Value = 0.07420673321644777
$ws.Range("F16").Value = 6.658088647926718
$ws.Range("J16").Value = 0.2512413524755317
$ws.Range("K16").Value = 0.8720427262133512
$ws.Range("M16").Value = 0.3004658737157655

$ws.Range("B17").Value = 0.8639935084187584
$ws.Range("C17").Value = 0.08668588425703661
$ws.Range("D17").Value = 0.02749995592343879
$ws.Range("E17").Value = 0.0744433328854246
$ws.Range("F17").Value = 6.560937960493249
$ws.Range("J17").Value = 0.2491488578229166
$ws.Range("K17").Value = 0.8629001139481431
$ws.Range("M17").Value = 0.2991269468428754

$ws.Range("B18").Value = 0.8584546077207165
$ws.Range("C18").Value = 0.08632657670430888
$ws.Range("D18").Value = 0.02701643123022279
$ws.Range("E18").Value = 0.07458448560307218
$ws.Range("F18").Value = 6.505119748026146
$ws.Range("J18").Value = 0.247951110029291
$ws.Range("K18").Value = 0.8577236255983678
$ws.Range("M18").Value = 0.2983839382441857

$ws.Range("B19").Value = 0.8565920435277405
$ws.Range("C19").Value = 0.08620697564580837
$ws.Range("D19").Value = 0.02685328807449849
$ws.Range("E19").Value = 0.07463314792722997
$ws.Range("F19").Value = 6.486230891248567
$ws.Range("J19").Value = 0.2475465671855233
$ws.Range("K19").Value = 0.8559850372587903
$ws.Range("M19").Value = 0.2981370203815459

$ws.Range("B20").Value = 0.8650247239912972
$ws.Range("C20").Value = 0.08675336087419794
$ws.Range("D20").Value = 0.02758971744964356
$ws.Range("E20").Value = 0.07441762209685443
$ws.Range("F20").Value = 6.571273547251565
$ws.Range("J20").Value = 0.2493710066814288
$ws.Range("K20").Value = 0.8638648590045932
$ws.Range("M20").Value = 0.2992666716171399

$ws.Range("B21").Value = 0.894250842006187
$ws.Range("C21").Value = 0.08873269661363281
$ws.Range("D21").Value = 0.0301044096556069
$ws.Range("E21").Value = 0.07375400794178688
$ws.Range("F21").Value = 6.857734866146529
$ws.Range("J21").Value = 0.2555703282154127
$ws.Range("K21").Value = 0.8913223342827621
$ws.Range("M21").Value = 0.3033862760384878

$ws.Range("B22").Value = 0.9140203154405242
$ws.Range("C22").Value = 0.09013407387287486
$ws.Range("D22").Value = 0.03177888230605674
$ws.Range("E22").Value = 0.07336537407935495
$ws.Range("F22").Value = 7.045518971013394
$ws.Range("J22").Value = 0.259674076810299
$ws.Range("K22").Value = 0.9100028477684816
$ws.Range("M22").Value = 0.3063208684494541

$ws.Range("B23").Value = 0.9034076365115027
$ws.Range("C23").Value = 0.08937622166136805
$ws.Range("D23").Value = 0.03088229421398125
$ws.Range("E23").Value = 0.07356867279157875
$ws.Range("F23").Value = 6.945241981487811
$ws.Range("J23").Value = 0.2574790348252236
$ws.Range("K23").Value = 0.8999651958801849
$ws.Range("M23").Value = 0.3047324204470243

$ws.Range("B24").Value = 0.864558286740845
$ws.Range("C24").Value = 0.08672281789117164
$ws.Range("D24").Value = 0.02754912653271901
$ws.Range("E24").Value = 0.07442922996948553
$ws.Range("F24").Value = 6.566600719721322
$ws.Range("J24").Value = 0.2492705568205196
$ws.Range("K24").Value = 0.8634284494901294
$ws.Range("M24").Value = 0.299203418706405

$ws.Range("B25").Value = 0.8255062694879314
$ws.Range("C25").Value = 0.08431838694248484
$ws.Range("D25").Value = 0.02408109951947068
$ws.Range("E25").Value = 0.07555122609501019
$ws.Range("F25").Value = 6.160462162629187
$ws.Range("J25").Value = 0.2406378098332027
$ws.Range("K25").Value = 0.8271538737769504
$ws.Range("M25").Value = 0.2942753584471518
